$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.5692650327263563; C = 0.5692650327263563; D = 0.3908711553892542; E = 0.6251968932978268; F = 0.2682304822149452; G = 14 }
    3  = @{ B = 0.3997456872501602; C = 0.4033467742491634; D = 0.2178465991588477; E = 0.4667403980360472; F = 0.2507737693766841; G = 13 }
    4  = @{ B = 0.3400043916411057; C = 0.3606492384221238; D = 0.1784350047575419; E = 0.4224156776891003; F = 0.261809129132386;  G = 12 }
    5  = @{ B = 0.4288908803047028; C = 0.4288908803047028; D = 0.2376241711832545; E = 0.4874670975391617; F = 0.242990663137872;  G = 11 }
    6  = @{ B = 0.3767221034172891; C = 0.3801506752190701; D = 0.1848819966348901; E = 0.429979065344919;  F = 0.2184858333361773; G = 10 }
    7  = @{ B = 0.3516807317407905; C = 0.358904738120446;  D = 0.1712159718543759; E = 0.4137825175794355; F = 0.2312546521125961; G = 9 }
    8  = @{ B = 0.3826977698461033; C = 0.3826977698461033; D = 0.193400430583437;  E = 0.4397731580979414; F = 0.2316224700134891; G = 8 }
    9  = @{ B = 0.3590813353308283; C = 0.3590813353308283; D = 0.1682323819789578; E = 0.4101614096657044; F = 0.2141069966830871; G = 7 }
    10 = @{ B = 0.3470367157308435; C = 0.3470367157308435; D = 0.1621328538698305; E = 0.4026572411739673; F = 0.2236918553848042; G = 6 }
    11 = @{ B = 0.308333697360539;  C = 0.308333697360539;  D = 0.1232078055746762; E = 0.3510096944169437; F = 0.1875437837101504; G = 5 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
